# Generate Report for Handback
# Updates the handback-status report after a new handback run:
#  - Overview sheet: status flips to "not in sync", refresh latest date
#  - zh-cn sheet: refresh handoff/handback timestamps
#  - de-de sheet: refresh handoff/handback timestamps
#  - Status column widened on all three sheets so the longer status text fits

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: now out of sync with en-US -------------------------------
$newStatus = "Handed back: not in sync with en-US"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value     = $newStatus
$dede.Range("C2").Value     = $newStatus

# --- Refreshed timestamps ----------------------------------------------------
$overview.Range("G2").Value = "2017-02-09 14:17:41"

$zhcn.Range("H2").Value = "2017-02-09 14:17:22"
$zhcn.Range("L2").Value = "2017-02-09 14:19:03"

$dede.Range("H2").Value = "2017-02-09 14:17:41"
$dede.Range("L2").Value = "2017-02-09 14:19:28"

# --- Widen the Status columns to fit the longer text -------------------------
$overview.Columns.Item(5).ColumnWidth = 32.65
$overview.Columns.Item(6).ColumnWidth = 32.65
$zhcn.Columns.Item(3).ColumnWidth = 32.65
$dede.Columns.Item(3).ColumnWidth = 32.65
